# Excess mortality - Week 50 update
# Applies the data corrections / additions described by the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
[void]$ws.Activate()

# --- Revised observed-mortality (column G) figures for several existing weeks ---
$ws.Range("G22").Value = 2673   # week 30
$ws.Range("G23").Value = 2668   # week 31
$ws.Range("G30").Value = 2719   # week 38
$ws.Range("G32").Value = 2997   # week 40
$ws.Range("G34").Value = 3216   # week 42
$ws.Range("G35").Value = 3445   # week 43
$ws.Range("G36").Value = 3675   # week 44
$ws.Range("G38").Value = 3560   # week 46
$ws.Range("G39").Value = 3317   # week 47
$ws.Range("G40").Value = 3388   # week 48
$ws.Range("G41").Value = 3494   # week 49

# --- New week 50 row ---
$ws.Range("F42").Value = 50
$ws.Range("G42").Value = 3571
$ws.Range("H42").Value = 3100
$ws.Range("I42").Formula = "=G42-H42"

# --- Move the "Som week 11 tot en met 19" totals row from row 43 down to row 46,
#     before row 43's own contents get overwritten with the week-51 number ---
$totalLabel = $ws.Range("F43").Value2
$ws.Range("F46").Value2 = $totalLabel

$ws.Range("G46").Formula = "=SUM(G3:G28)"
$ws.Range("H46").Formula = "=SUM(H3:H28)"
$ws.Range("I46").Formula = "=SUM(I3:I34)"
$ws.Range("G46").NumberFormat = $ws.Range("G43").NumberFormat
$ws.Range("H46").NumberFormat = $ws.Range("H43").NumberFormat
$ws.Range("I46").NumberFormat = $ws.Range("I43").NumberFormat

# Wipe the old totals row (contents + formatting) - it's replaced by row 46
[void]$ws.Range("F43:I43").Clear()

# --- Week numbers 51 and 52, currently without observed/expected data yet ---
$ws.Range("F43").Value = 51
$ws.Range("F44").Value = 52

# --- Update the view: scroll position + active cell selection ---
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
[void]$ws.Range("I43").Select()

Write-Host "Week 50 update applied"
